$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 16 de Mayo de 2020 a las 13:35"

# Update the data rows (Casos totales, Casos activos, Recuperados, Muertes)
# Row 4
$ws.Range("B4").Value = 66210
$ws.Range("C4").Value = 40586
$ws.Range("D4").Value = 16798
$ws.Range("E4").Value = 8826

# Row 5
$ws.Range("B5").Value = 55685
$ws.Range("C5").Value = 26067
$ws.Range("D5").Value = 23703
$ws.Range("E5").Value = 5915

# Row 6
$ws.Range("B6").Value = 18369
$ws.Range("C6").Value = 7735
$ws.Range("D6").Value = 8694
$ws.Range("E6").Value = 1940

# Row 7
$ws.Range("B7").Value = 16587
$ws.Range("C7").Value = 6328
$ws.Range("D7").Value = 7376
$ws.Range("E7").Value = 2883

# Row 9
$ws.Range("B9").Value = 12420
$ws.Range("C9").Value = 10319
$ws.Range("D9").Value = 746
$ws.Range("E9").Value = 1355

# Row 11
$ws.Range("B11").Value = 9041
$ws.Range("C11").Value = 8409
$ws.Range("D11").Value = 28
$ws.Range("E11").Value = 604

# Row 14
$ws.Range("B14").Value = 5456
$ws.Range("C14").Value = 3655
$ws.Range("D14").Value = 963
$ws.Range("E14").Value = 838

# Row 16
$ws.Range("B16").Value = 5137
$ws.Range("C16").Value = 3658
$ws.Range("D16").Value = 978
$ws.Range("E16").Value = 501

# Row 20
$ws.Range("B20").Value = 4020
$ws.Range("C20").Value = 2990
$ws.Range("D20").Value = 682

# Row 32
$ws.Range("B32").Value = 2366
$ws.Range("C32").Value = 1055
$ws.Range("D32").Value = 998
$ws.Range("E32").Value = 313

# Row 33
$ws.Range("B33").Value = 2284
$ws.Range("C33").Value = 1524
$ws.Range("D33").Value = 609
